# Daily attendance processing - 2026-01-03 11:30:58
# Swap the order of the "System" and email address tokens in column G
# wherever the cell currently reads "System, dnasr281@gmail.com" so it
# becomes "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
